$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 44207
$ws.Range("B2").Value = 44216
$ws.Range("C2").Value = 727695.1500000001
$ws.Range("D2").Value = 88097.13499999998
$ws.Range("E2").Value = 496711.2499999999
$ws.Range("F2").Value = 403027.01
$ws.Range("G2").Value = 675741.05
$ws.Range("H2").Value = 131857.4809999999
$ws.Range("I2").Value = 0.07139541880964863
